$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.783.80'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '3.383.19'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'580.24"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = "'177.99"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E7').Value = '  +4.12%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '3.384.91'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '3.972.97'
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = "'28.84"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.44%  '
$ws.Range('D16').Value = '65.938.83'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '3.381.95'
$ws.Range('E18').Value = '  -1.92%  '
$ws.Range('D19').Value = "'5.85"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = "'13.69"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'365.08"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').Value = "'72.56"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = "'0.997"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  +4.50%  '
$ws.Range('D27').Value = "'9.73"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = "'5.72"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').Value = "'23.05"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('D33').Value = "'1.00"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = "'6.97"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('D35').Value = "'1.25"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.09%  '
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').Value = "'162.48"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('D38').Value = "'0.858"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.35%  '
$ws.Range('D39').Value = "'27.19"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.68%  '
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').Value = "'2.60"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').Value = '2.675.76'
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('D43').Value = "'4.34"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').Value = "'6.19"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = "'0.0680"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = "'24.50"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').Value = "'39.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D48').Value = "'330.97"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.95%  '
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D51').Value = "'31.36"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.63%  '
